{"js": "// 1) \"Our problem domain includes ...\" -> \"The problem domain includes ...\"\nconst ourResults = context.document.body.search(\"Our problem domain includes\", { matchCase: true, matchWholeWord: false });\nourResults.load(\"text\");\nawait context.sync();\nif (ourResults.items.length > 0) {\n  ourResults.items[0].insertText(\"The problem domain includes\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"Testing should be silent, automatic and continuous\" -> \"Testing should be automatic and continuous\"\nconst testingResults = context.document.body.search(\"Testing should be silent, automatic and continuous\", { matchCase: true, matchWholeWord: false });\ntestingResults.load(\"text\");\nawait context.sync();\nif (testingResults.items.length > 0) {\n  testingResults.items[0].insertText(\"Testing should be automatic and continuous\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Append \"   9454948223\" after \"Incremental\" in its own run/paragraph content\nconst incResults = context.document.body.search(\"Incremental\", { matchCase: true, matchWholeWord: true });\nincResults.load(\"text\");\nawait context.sync();\nif (incResults.items.length > 0) {\n  incResults.items[0].insertText(\"   9454948223\", \"After\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Our problem domain includes\" -> \"The problem domain includes\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Text = \"Our problem domain includes\"\n$rng1.Find.MatchCase = $true\n$rng1.Find.MatchWholeWord = $false\n$rng1.Find.Execute(\"Our problem domain includes\", $false, $false, $false, $false, $false, $true, 1, $false, \"The problem domain includes\", 2)\n\n# 2) \"Testing should be silent, automatic and continuous\" -> \"Testing should be automatic and continuous\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"Testing should be silent, automatic and continuous\"\n$rng2.Find.MatchCase = $true\n$rng2.Find.MatchWholeWord = $false\n$rng2.Find.Execute(\"Testing should be silent, automatic and continuous\", $false, $false, $false, $false, $false, $true, 1, $false, \"Testing should be automatic and continuous\", 2)\n\n# 3) Append \"   9454948223\" right after the standalone \"Incremental\" paragraph text\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Text = \"Incremental\"\n$rng3.Find.MatchCase = $true\n$rng3.Find.MatchWholeWord = $true\n$found3 = $rng3.Find.Execute()\nif ($found3) {\n    $rng3.InsertAfter(\"   9454948223\")\n}\n"}
